$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns D, M, N, O, P, S
$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr3 = "${col}3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}
